$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poses")

# --- Remove stray "TODO: Re-derive X Y Z" note from H1 (its shared string
#     becomes unused and is pruned on save, which also renumbers every
#     later shared-string index used elsewhere on this sheet) ---
$ws.Range("H1").ClearContents()

# --- New column A width + selection moved to J27 ---
$ws.Columns.Item(1).ColumnWidth = 15.498697916666666
$ws.Range("J27").Select()

# --- Joint offset fixes: X offsets recentred to 0, Z offsets moved to the
#     midpoint of their segment (instead of the far end) for the ankle/
#     knee/hip/spine/elbow joint rows. Columns E/F/G are downstream
#     shared formulas and recalc automatically. ---

# left_ankle / right_ankle (row 3 / row 5): X -> 0
$ws.Range("C3").Value = 0
$ws.Range("C5").Value = 0

# left_knee (row 7): X -> 0, Z -> Height!E6/2
$ws.Range("C7").Value = 0
$ws.Range("D7").Formula = '=Height!$E$6/2'

# right_knee (row 9): X -> 0, Z -> Height!E6/2
$ws.Range("C9").Value = 0
$ws.Range("D9").Formula = '=Height!$E$6/2'

# left_hip (row 11): X -> 0, Z -> Height!E7/2
$ws.Range("C11").Value = 0
$ws.Range("D11").Formula = '=Height!$E$7/2'

# right_hip (row 13): X -> 0, Z -> Height!E7/2
$ws.Range("C13").Value = 0
$ws.Range("D13").Formula = '=Height!$E$7/2'

# lower_spine (row 15): Z -> 0
$ws.Range("D15").Value = 0

# upper_spine (row 17): Z -> 0
$ws.Range("D17").Value = 0

# neck (row 19): Z -> -Height!E8/2
$ws.Range("D19").Formula = '=-Height!$E$8/2'

# left_shoulder (row 22): X -> 0, Z -> -Height!E4/2
$ws.Range("C22").Value = 0
$ws.Range("D22").Formula = '=-Height!$E$4/2'

# right_shoulder (row 24): X -> 0, Z -> -Height!E4/2
$ws.Range("C24").Value = 0
$ws.Range("D24").Formula = '=-Height!$E$4/2'

# left_elbow (row 26): X -> 0, Z -> -Height!E3/2
$ws.Range("C26").Value = 0
$ws.Range("D26").Formula = '=-Height!$E$3/2'

# right_elbow (row 28): X -> 0, Z -> -Height!E3/2
$ws.Range("C28").Value = 0
$ws.Range("D28").Formula = '=-Height!$E$3/2'

# left_wrist (row 30): X -> 0, Z -> -Height!E2/2
$ws.Range("C30").Value = 0
$ws.Range("D30").Formula = '=-Height!$E$2/2'

# right_wrist (row 32): X -> 0, Z -> -Height!E2/2
$ws.Range("C32").Value = 0
$ws.Range("D32").Formula = '=-Height!$E$2/2'
